$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 2017 "GL_alvar" row (previously row 8) up into the
# gap at row 7, and re-point the two rows below it (MB_alvar / Prairie)
# so the 2017 block becomes contiguous rows 7-9.
$ws.Range("A7").Value = 2017
$ws.Range("B7").Value = "GL_alvar"
$ws.Range("C7").Value = 544.65099999999995
$ws.Range("D7").Value = 163.02699999999999

$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = "MB_alvar"
$ws.Range("C8").Value = 93.6
$ws.Range("D8").Value = 81.897999999999996

$ws.Range("A9").Value = 2017
$ws.Range("B9").Value = "Prairie"
$ws.Range("C9").Value = 70.379000000000005
$ws.Range("D9").Value = 75.872

# Add the newly completed 2018 region-specific fitness estimates.
$ws.Range("A10").Value = 2018
$ws.Range("B10").Value = "GL_alvar"
$ws.Range("C10").Value = 926.32399999999996
$ws.Range("D10").Value = 268.49599999999998

$ws.Range("A11").Value = 2018
$ws.Range("B11").Value = "MB_alvar"
$ws.Range("C11").Value = 248.6
$ws.Range("D11").Value = 194.65

$ws.Range("A12").Value = 2018
$ws.Range("B12").Value = "Prairie"
$ws.Range("C12").Value = 94.444999999999993
$ws.Range("D12").Value = 100.45399999999999

# Column D needed a bit more room once the new values were added.
$ws.Columns.Item(4).ColumnWidth = 10.0

# Reflect the new active cell/selection at the bottom of the table.
$ws.Range("D12").Select()
